$wb = $excel.ActiveWorkbook

# --- Sheet "disk_list": swap row 1 and row 2 contents ---
$ws = $wb.Worksheets.Item("disk_list")

# Column A holds device codes stored as text ("514-2", "518"). "518" looks
# numeric, so a plain Value assignment would get auto-converted to a number.
# Use Cut (via an unused scratch cell) to swap A1/A2 while preserving their
# original text type and not touching any cell styling.
$ws.Range("A1").Cut($ws.Range("H1")) | Out-Null
$ws.Range("A2").Cut($ws.Range("A1")) | Out-Null
$ws.Range("H1").Cut($ws.Range("A2")) | Out-Null

# Columns B..E: read the old values (Value2; the .Value getter is unreliable
# in this runtime) and write them back swapped by position. This keeps any
# position-anchored formatting (e.g. C1's Hyperlink style) exactly where it
# was, matching how the style stayed on C1 across the edit.
$b1 = $ws.Range("B1").Value2
$c1 = $ws.Range("C1").Value2
$d1 = $ws.Range("D1").Value2
$e1 = $ws.Range("E1").Value2

$b2 = $ws.Range("B2").Value2
$c2 = $ws.Range("C2").Value2
$d2 = $ws.Range("D2").Value2
$e2 = $ws.Range("E2").Value2

$ws.Range("B1").Value = $b2
$ws.Range("C1").Value = $c2
$ws.Range("D1").Value = $d2
$ws.Range("E1").Value = $e2
$ws.Range("F1").Value = "první sít, ixon"

$ws.Range("B2").Value = $b1
$ws.Range("C2").Value = $c1
$ws.Range("D2").Value = $d1
$ws.Range("E2").Value = $e1
$ws.Range("F2").ClearContents()

# --- Sheet "Settings": update B3, B4, B5 ---
$ws2 = $wb.Worksheets.Item("Settings")
$ws2.Range("B3").Value = 1
$ws2.Range("B4").Value = 1
$ws2.Range("B5").Value = 0
